# Rebuild the single placeholder paragraph as the FRA template "notes"
# block: a RE:/Date of Crash/Date of Birth header followed by the
# doc_body placeholder paragraph. We insert the replacement content as
# raw WordprocessingML so the exact run/tab/proofErr structure produced
# by Word is reproduced faithfully.

$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$r = $p.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">RE: </w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>case_name</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'

$para2 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' + `
    '</w:p>'

$para3 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' + `
    '<w:r><w:t>Date of Crash:</w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>crash_date</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'

$para4 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' + `
    '<w:r><w:t>Date of Birth:</w:t></w:r>' + `
    '<w:r><w:tab/><w:t>plaintiff1_</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>dob</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'

$para5 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' + `
    '<w:r><w:tab/><w:t>plaintiff2_</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>dob</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'

$para6 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="2127" w:hanging="2127"/></w:pPr>' + `
    '</w:p>'

$para7 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:spacing w:after="0"/><w:ind w:left="2126" w:hanging="2126"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>doc_body</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'

$xml = $para1 + $para2 + $para3 + $para4 + $para5 + $para6 + $para7

[void]$r.InsertXML($xml)
